$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Header values ---
# G1 used to be "Initial Comment"; it now becomes "Category IDs".
$ws.Range("G1").Value = "Category IDs"
# H1 becomes the (relocated) "Initial Comment" header.
$ws.Range("H1").Value = "Initial Comment"
# I1 is a brand new "Additional Comments" header.
$ws.Range("I1").Value = "Additional Comments"

# --- Header styles (grey = required, blue = optional) ---
# G1 moves from required(grey) to optional(blue) - copy format from F1 (Priority, optional).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
# H1 takes over the required(grey) look the old "Initial Comment" cell (G1) used to have - copy from A1 (Subject, required).
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
# I1 is optional(blue) - copy from F1 too.
$ws.Range("F1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Comments ---
# G1's comment changes from the "initial comment" help text to the new "category ids" help text.
$ws.Range("G1").Comment.Text("Insert a list of category IDs here, separated by comma.")
# H1 gets a new comment carrying over the old "initial comment" help text.
$ws.Range("H1").AddComment("An initial comment to open the ticket with. ")
# I1 gets a new comment explaining additional comment columns.
$ws.Range("I1").AddComment("You can continue to add additional comments as additional columns, including and after this one. This allows you to insert as many comments as you want on a ticket.")

Write-Host "done"
